$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, centered, bordered) onto the new
# header cells I1:J1 before setting their text, so they match the rest
# of the header row's formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data rows: row number, I value, J value
$data = @(
    @(2,6,6),
    @(3,7,7),
    @(4,7,7),
    @(5,8,8),
    @(6,9,9),
    @(7,10,10),
    @(8,9,9),
    @(9,7,7),
    @(10,8,8),
    @(11,8,8),
    @(12,8,8),
    @(13,8,8),
    @(14,7,7),
    @(15,9,9),
    @(16,6,6),
    @(17,7,8),
    @(18,9,9),
    @(19,9,9),
    @(20,9,9),
    @(21,9,9),
    @(22,9,9),
    @(23,9,9),
    @(24,9,9),
    @(25,9,9),
    @(26,9,9),
    @(27,9,9),
    @(28,9,9),
    @(29,9,9),
    @(30,9,9),
    @(31,9,9),
    @(32,9,9),
    @(33,9,9),
    @(34,9,9),
    @(35,9,9),
    @(36,9,9),
    @(37,9,9),
    @(38,9,9),
    @(39,9,9),
    @(40,9,9),
    @(41,9,9),
    @(42,10,10),
    @(43,9,9),
    @(44,10,10),
    @(45,9,9),
    @(46,9,9),
    @(47,9,9),
    @(48,9,9),
    @(49,8,9),
    @(50,9,9),
    @(51,9,9),
    @(52,9,9),
    @(53,9,9),
    @(54,9,9),
    @(55,9,9),
    @(56,9,9),
    @(57,9,9),
    @(58,9,9),
    @(59,9,9),
    @(60,8,9),
    @(61,9,9),
    @(62,9,9),
    @(63,9,9),
    @(64,6,6),
    @(65,5,5),
    @(66,4,4),
    @(67,4,4),
    @(68,5,5),
    @(69,4,4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
